# memperbaiki struktur table user, dari prodi menjadi tingkat
# Adds a new "Tingkat (Opsional)" column (H) to the user table and fills
# in a value for each row, plus the previously-missing "Angkatan" value
# for the last row (G5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlLeft = -4131

# New column header (H1) - same look as the rest of the header row
# (bold font, centered horizontally & vertically).
$ws.Range("H1").Value = "Tingkat (Opsional)"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = $xlCenter
$ws.Range("H1").VerticalAlignment = $xlCenter

# New column values per row (H2:H5) - same look as the other data cells
# (regular font, left aligned, vertically centered).
$ws.Range("H2").Value = "Tingkat IV"
$ws.Range("H3").Value = "Tingkat III"
$ws.Range("H4").Value = "Tingkat III"
$ws.Range("H5").Value = "Tingkat I"
$ws.Range("H2:H5").HorizontalAlignment = $xlLeft
$ws.Range("H2:H5").VerticalAlignment = $xlCenter

# Fill in the previously-missing Angkatan value for row 5
$ws.Range("G5").Value = 64

# Resize column widths: narrow column B, add a width for the new column H
$ws.Columns.Item(2).ColumnWidth = 26.8
$ws.Columns.Item(8).ColumnWidth = 17.65

# Update the active selection to mirror what Excel would show after editing
$ws.Range("G9").Select() | Out-Null
